# Updates the last OHLCV row (row 132) with revised close/volume figures and
# appends two new daily rows (133, 134) pulled in from the latest Binance
# ETHUSDT data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping it stored as TEXT
# (this sheet keeps price/volume columns as inline strings, not numbers),
# and keep the cell's style at the workbook default (no borders/bold/etc).
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Helper: write a true numeric value.
function Set-NumCell($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# Insert two fresh rows after the current last row (132) for the new
# 2020-06-11 and 2020-06-12 daily candles.
$ws.Range("A133:A134").EntireRow.Insert()

# Carry the bordered/bold "index" column style from row 132 down into the
# two newly inserted rows so column A keeps its look.
$ws.Cells.Item(132, 1).Copy()
$ws.Range("A133:A134").PasteSpecial(-4122)

# Row 132 updates
Set-TextCell 132 3 "250.28000000"
Set-TextCell 132 4 "242.00000000"
Set-TextCell 132 5 "247.78000000"
Set-TextCell 132 6 "431285.61715000"
Set-TextCell 132 8 "105797248.24742570"
Set-NumCell 132 9 185330
Set-TextCell 132 10 "221305.11120000"
Set-TextCell 132 11 "54322143.50039510"
Set-NumCell 132 13 247.7799999999998
Set-NumCell 132 14 245.7899999999999
Set-NumCell 132 15 243.96
Set-NumCell 132 16 237.5053333333335
Set-NumCell 132 17 220.4876666666666
Set-NumCell 132 18 247.78
Set-NumCell 132 19 246.6380655797366
Set-NumCell 132 20 239.2247235404768
Set-NumCell 132 21 227.662815602148
Set-NumCell 132 22 11.56190793832886
Set-NumCell 132 23 11.00574053443424
Set-NumCell 132 24 0.5561674038946158

# Row 133 (new)
Set-NumCell 133 1 131
Set-TextCell 133 2 "247.78000000"
Set-TextCell 133 3 "250.09000000"
Set-TextCell 133 4 "226.20000000"
Set-TextCell 133 5 "230.51000000"
Set-TextCell 133 6 "816456.89161000"
Set-NumCell 133 7 1591919999999
Set-TextCell 133 8 "195063781.27644910"
Set-NumCell 133 9 312655
Set-TextCell 133 10 "373492.34304000"
Set-TextCell 133 11 "89279896.48515260"
Set-TextCell 133 12 "2020-06-11 08:00:00"
Set-NumCell 133 13 230.5099999999998
Set-NumCell 133 14 239.1449999999999
Set-NumCell 133 15 242.1457142857143
Set-NumCell 133 16 238.9866666666668
Set-NumCell 133 17 221.846
Set-NumCell 133 18 230.51
Set-NumCell 133 19 235.8860218599122
Set-NumCell 133 20 237.8839968415866
Set-NumCell 133 21 227.8737263190247
Set-NumCell 133 22 10.01027052256188
Set-NumCell 133 23 10.80664653205974
Set-NumCell 133 24 -0.7963760094978554

# Row 134 (new)
Set-NumCell 134 1 132
Set-TextCell 134 2 "230.46000000"
Set-TextCell 134 3 "239.38000000"
Set-TextCell 134 4 "228.19000000"
Set-TextCell 134 5 "235.22000000"
Set-TextCell 134 6 "397549.31654000"
Set-NumCell 134 7 1592006399999
Set-TextCell 134 8 "93464813.52779320"
Set-NumCell 134 9 154895
Set-TextCell 134 10 "199267.27950000"
Set-TextCell 134 11 "46876827.74583840"
Set-TextCell 134 12 "2020-06-12 08:00:00"
Set-NumCell 134 13 235.2199999999998
Set-NumCell 134 14 232.8649999999999
Set-NumCell 134 15 241.4585714285714
Set-NumCell 134 16 239.9873333333335
Set-NumCell 134 17 223.034
Set-NumCell 134 18 235.22
Set-NumCell 134 19 235.4420072866374
Set-NumCell 134 20 237.4741511735583
Set-NumCell 134 21 228.4179142557757
Set-NumCell 134 22 9.05623691778257
Set-NumCell 134 23 10.45656460920426
Set-NumCell 134 24 -1.40032769142169

Write-Output "done"
